$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for rows 2 through 28 from 45482 to 45483
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45483
}

# Ensure row 28 has an explicit row height (15, custom)
$ws.Rows.Item(28).RowHeight = 15

# Add the new row 29 with its data
$ws.Cells.Item(29, 1).Value = "A 29276-2024"
$ws.Cells.Item(29, 2).Value = 45483
$ws.Cells.Item(29, 3).Value = 45483
$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"
$ws.Cells.Item(29, 7).Value = 1.8
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0

# Copy style from row 28 for B/C (date) and R (wrap text) columns
$ws.Cells.Item(28, 2).Copy()
$ws.Cells.Item(29, 2).PasteSpecial(-4122)
$ws.Cells.Item(28, 3).Copy()
$ws.Cells.Item(29, 3).PasteSpecial(-4122)
$ws.Cells.Item(28, 18).Copy()
$ws.Cells.Item(29, 18).PasteSpecial(-4122)
